$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.018±0.002"
$ws.Range("C2").Value = "0.186±0.013"

$ws.Range("B3").Value = "0.066±0.008"
$ws.Range("C3").Value = "0.245±0.010"

$ws.Range("B4").Value = "0.899±0.033"
$ws.Range("C4").Value = "0.480±0.037"

$ws.Range("B5").Value = "0.994±0.001"
$ws.Range("C5").Value = "0.416±0.043"

$ws.Range("B6").Value = "0.985±0.009"
$ws.Range("C6").Value = "0.800±0.010"

$ws.Range("B7").Value = "0.962±0.021"
$ws.Range("C7").Value = "0.330±0.031"

$ws.Range("B8").Value = "0.007±0.004"
$ws.Range("C8").Value = "0.177±0.005"
